$p = $ppt.ActivePresentation

# Slides 3 through 13 each contain an empty, unused placeholder shape
# named "Content Placeholder 2" (or "Text Placeholder 2" on the final
# "Thank You" slide) that the template generator left behind with no
# content. Remove each of these empty placeholder shapes.
for ($i = 3; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = $s.Shapes.Count; $j -ge 1; $j--) {
        $shp = $s.Shapes.Item($j)
        $name = $shp.Name
        if ($name -eq "Content Placeholder 2" -or $name -eq "Text Placeholder 2") {
            $shp.Delete()
        }
    }
}
